# Insert a new data row at row 184 (pushing existing rows 184-203 down to 185-204)
# and populate it with the new weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 184; this shifts rows 184-203 down to 185-204
$ws.Rows.Item(184).Insert()

# Populate the new row 184 with the new record
$ws.Range("A184").Value = 5
$ws.Range("B184").Value = "Macroferia Regional de Talca"
$ws.Range("C184").Value = "Maule"
$ws.Range("D184").Value = 45218
$ws.Range("E184").Value = 7
$ws.Range("F184").Value = "Fruta"
$ws.Range("G184").Value = 100108
$ws.Range("H184").Value = "Tropicales y subtropicales"
$ws.Range("I184").Value = 100108002
$ws.Range("J184").Value = "Mango"
$ws.Range("K184").Value = "Sin especificar"
$ws.Range("L184").Value = "Primera"
$ws.Range("M184").Value = 248
$ws.Range("N184").Value = 9000
$ws.Range("O184").Value = 9000
$ws.Range("P184").Value = 9000
$ws.Range("Q184").Value = "`$/bandeja 4 kilos"
$ws.Range("R184").Value = "Brasil"
$ws.Range("S184").Value = 2250
$ws.Range("T184").Value = 4
